{"js": "// Add a new paragraph \"Test!\" right after the existing last paragraph\n// (i.e. before the final section properties), matching the diff which\n// inserts a new <w:p><w:r><w:t>Test!</w:t></w:r></w:p> after the\n// \"Your resume or CV file \" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Test!\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add a new paragraph \"Test!\" at the end of the document, right after the\n# existing \"Your resume or CV file \" paragraph.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Collapse(0)            # wdCollapseEnd\n$range.InsertParagraphAfter()\n$range.Collapse(0)            # wdCollapseEnd\n$range.Text = \"Test!\"\n"}
